# Refresh cryptocurrency price/volume snapshot (scheduled GitHub Actions data pull)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '63.906.52'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.53%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.138.23'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.45%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '587.35'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.02%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '146.43'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.67%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.136.23'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.28%  '
$ws.Range('E9').Value = '  -1.05%  '
$ws.Range('E10').Value = '  -0.14%  '
$ws.Range('E11').Value = '  -0.01%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.459'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.71%  '
$ws.Range('E13').Value = '  -3.17%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '36.93'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.41%  '
$ws.Range('E15').Value = '  -1.68%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.658.35'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.42%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '63.711.25'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.69%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.135.31'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.54%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.08'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.36%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '465.59'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.00%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.30'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.05%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.734'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.13%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.44'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.46%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.93'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -3.60%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '81.22'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.67%  '
$ws.Range('E26').Value = '  +0.23%  '
$ws.Range('E27').Value = '  -0.12%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.34'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +7.30%  '
$ws.Range('B29').Value = 'ImmutableX'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.23'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.01%  '
$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.68'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.23%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.00'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.07%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.05'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.48%  '
$ws.Range('E33').Value = '  -1.07%  '
$ws.Range('E34').Value = '  +0.27%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0₃0850'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.35%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.04'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.08%  '
$ws.Range('E37').Value = '  -5.22%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.32'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.77%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '51.19'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.31%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '438.71'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.59%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.87'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.16%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.922.80'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.67%  '
$ws.Range('E44').Value = '  -0.44%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.280'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.65%  '
$ws.Range('E46').Value = '  -5.43%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '36.77'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +4.90%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '127.21'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +3.15%  '
$ws.Range('E50').Value = '  -1.46%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '24.13'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.27%  '
